# Auto-update draw results: append the 2025-10-30 Pick 4 draw as a new
# row (44) at the bottom of the "Results" sheet.
#
# Columns: A=Date, B=Game, C=Phase, D=Result, E=InsertedAt
# All five columns hold plain text in the source data (e.g. "251030" is
# kept as text, not parsed as a number; "2025-10-30" is kept as text,
# not parsed as a date) so we pre-format the numeric-looking columns
# (A and C) as Text before writing into them - this mirrors how Excel
# avoids auto-converting typed text that looks like a number/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-10-30"

$ws.Range("B$row").Value = "Pick 4"

$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = "251030"

$ws.Range("D$row").Value = "6-9-8-7"

$ws.Range("E$row").Value = "2025-10-30T21:40:12.400+04:00"
